# Updates the cryptos list (price/volume columns) with refreshed figures from the
# Sat Nov  2 17:57:02 UTC 2024 GitHub Actions data pull. Also reflects the two
# ranking swaps captured in this snapshot: Kaspa/WhiteBITCoin (rows 37-38) and
# Cronos/POPCAT (rows 49-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # The Price column stores numeric-looking strings (e.g. "568.77") as plain
    # text in the source data. Force text format before assigning so Excel does
    # not reinterpret the string as a number, then drop the temporary format so
    # no stray cell style is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "69.384.39"
$ws.Range("E2").Value = "  -0.61%  "
Set-TextValue "D3" "2.484.86"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "568.77"
$ws.Range("E5").Value = "  -1.35%  "
Set-TextValue "D6" "164.06"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.12%  "
Set-TextValue "D9" "2.482.35"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  -1.11%  "
Set-TextValue "D14" "2.942.77"
$ws.Range("E14").Value = "  -1.21%  "
Set-TextValue "D15" "69.220.46"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  -1.27%  "
Set-TextValue "D17" "24.15"
$ws.Range("E17").Value = "  -3.40%  "
Set-TextValue "D18" "2.496.39"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -2.56%  "
Set-TextValue "D20" "7.35"
$ws.Range("E20").Value = "  -4.88%  "
Set-TextValue "D21" "346.12"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -5.77%  "
$ws.Range("E24").Value = "  -0.14%  "
Set-TextValue "D25" "69.47"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("E26").Value = "  -3.30%  "
Set-TextValue "D28" "8.59"
$ws.Range("E28").Value = "  -4.37%  "
$ws.Range("E29").Value = "  +0.00%  "
Set-TextValue "D30" "0.0₃0864"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("E31").Value = "  -4.46%  "
Set-TextValue "D32" "435.47"
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D37" "0.112"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D38" "19.06"
$ws.Range("E38").Value = "  +0.12%  "
Set-TextValue "D39" "18.11"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("E41").Value = "  -2.43%  "
Set-TextValue "D42" "4.56"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("E43").Value = "  -2.78%  "
Set-TextValue "D44" "2.15"
$ws.Range("E44").Value = "  -6.31%  "
$ws.Range("E45").Value = "  -6.74%  "
Set-TextValue "D46" "137.57"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("B49").Value = "POPCAT"
$ws.Range("C49").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
Set-TextValue "D49" "2.02"
$ws.Range("E49").Value = "  +24.28%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0723"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -1.18%  "
